# 🤖 自动更新价格数据 2025-12-23 02:48:22
# Insert a new top row of price data (row 2), pushing all existing
# date rows down by one, and populate it with today's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the date column as plain text (it already holds literal
# "yyyy-mm-dd" strings, not real dates) so re-writing the shifted
# values below doesn't get auto-coerced into date serials.
$ws.Range("A2:A34").NumberFormat = "@"

# Existing data occupies rows 2..33 (below the header in row 1).
# Shift it down by one row, from the bottom up, copying values only
# (Value2, not Value, to dodge locale/date coercion on read) so no
# formatting/style gets pulled along for the ride.
for ($r = 33; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value2 = $ws.Cells.Item($src, 1).Value2
    $ws.Cells.Item($dst, 2).Value2 = $ws.Cells.Item($src, 2).Value2
    $ws.Cells.Item($dst, 3).Value2 = $ws.Cells.Item($src, 3).Value2
    $ws.Cells.Item($dst, 4).Value2 = $ws.Cells.Item($src, 4).Value2
}

# Populate the freshly vacated top data row with today's reading.
$ws.Cells.Item(2, 1).Value2 = "2025-12-23"
$ws.Cells.Item(2, 2).Value2 = 783.5
$ws.Cells.Item(2, 3).Value2 = 1112
$ws.Cells.Item(2, 4).Value2 = 3610
